$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "normal_tissue": update a few existing EGFR/skin rows and append new
# EGFR/skin cell-type rows (basal/corneal/spinous layers, endothelial cells,
# extracellular matrix, fibrohistiocytic cells, langerhans cells, lymphocytes,
# melanocytes, vascular mural cells).
# ---------------------------------------------------------------------------
$wsNT = $wb.Worksheets.Item("normal_tissue")

# breast / EGFR rows gain a "Low" reliability level
$wsNT.Range("E12").Value = "Low"
$wsNT.Range("E13").Value = "Low"
$wsNT.Range("E14").Value = "Low"

# skin 1 / EGFR rows: cell types get corrected + a couple of level updates
$wsNT.Range("D18").Value = "cells in basal layer"
$wsNT.Range("E18").Value = "Medium"
$wsNT.Range("D19").Value = "cells in corneal layer"
$wsNT.Range("D20").Value = "cells in granular layer"
$wsNT.Range("D21").Value = "cells in spinous layer"

# new rows appended after row 21
$newRows21 = @(
    @("ENSG00000146648", "EGFR", "skin 1", "endothelial cells",      "Not detected", "Enhanced"),
    @("ENSG00000146648", "EGFR", "skin 1", "extracellular matrix",   "Not detected", "Enhanced"),
    @("ENSG00000146648", "EGFR", "skin 1", "fibrohistiocytic cells", "High",         "Enhanced"),
    @("ENSG00000146648", "EGFR", "skin 1", "langerhans cells",       "Not detected", "Enhanced"),
    @("ENSG00000146648", "EGFR", "skin 1", "lymphocytes",            "Not detected", "Enhanced"),
    @("ENSG00000146648", "EGFR", "skin 1", "melanocytes",            "Not detected", "Enhanced"),
    @("ENSG00000146648", "EGFR", "skin 1", "vascular mural cells",   "Not detected", "Enhanced")
)

$r = 22
foreach ($row in $newRows21) {
    for ($i = 0; $i -lt $row.Length; $i++) {
        $wsNT.Cells.Item($r, $i + 1).Value = $row[$i]
    }
    $r++
}

# ---------------------------------------------------------------------------
# Sheet "pathology": refresh count columns and store the p-value-ish numbers
# as literal text (matches the "2.289e-2" style strings in the source data).
# ---------------------------------------------------------------------------
$wsPath = $wb.Worksheets.Item("pathology")

$wsPath.Range("K2").Value = "'2.289e-2"
$wsPath.Range("I3").Value = "'9.953e-2"
$wsPath.Range("K4").Value = "'2.750e-2"

$wsPath.Range("E5").Value = 1
$wsPath.Range("F5").Value = 2
$wsPath.Range("G5").Value = 9
$wsPath.Range("I5").Value = "'2.618e-2"

$wsPath.Range("D6").Value = 6
$wsPath.Range("E6").Value = 3
$wsPath.Range("G6").Value = 1
$wsPath.Range("I6").Value = "'1.079e-1"

$wsPath.Range("E7").Value = 4
$wsPath.Range("F7").Value = 1
$wsPath.Range("G7").Value = 6
$wsPath.Range("K7").Value = "'2.846e-2"

# ---------------------------------------------------------------------------
# Sheet "subcellular_location": TP53 row reclassified from single
# "Nucleoplasm" to a multi-location "Cytosol;Mitochondria;Nucleoplasm" call.
# ---------------------------------------------------------------------------
$wsSub = $wb.Worksheets.Item("subcellular_location")

$wsSub.Range("C2").Value = "Supported"
$wsSub.Range("E2").Value = "Cytosol;Mitochondria"
$wsSub.Range("G2").Value = ""
$wsSub.Range("H2").Value = "Cytosol;Mitochondria;Nucleoplasm"
$wsSub.Range("N2").Value = "Cytosol (GO:0005829);Mitochondria (GO:0005739);Nucleoplasm (GO:0005654)"
